$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text cells keep their exact string representation (no numeric
# auto-conversion / rounding / scientific notation) by forcing Text format
# before assigning price values in column D.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.200.50"
$ws.Range("E2").Value = "  +0.39%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.905.82"
$ws.Range("E3").Value = "  +0.80%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.88"
$ws.Range("E5").Value = "  +0.35%  "
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5207"
$ws.Range("E7").Value = "  +1.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3766"
$ws.Range("E8").Value = "  +0.21%  "
$ws.Range("E9").Value = "  +0.97%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.20"
$ws.Range("E10").Value = "  +0.35%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9045"
$ws.Range("E11").Value = "  +0.28%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08463"
$ws.Range("E12").Value = "  +10.45%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "96.93"
$ws.Range("E13").Value = "  +2.75%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.903.93"
$ws.Range("E14").Value = "  +0.64%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.299"
$ws.Range("E15").Value = "  +1.04%  "
$ws.Range("E16").Value = "  +0.11%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008648"
$ws.Range("E17").Value = "  +1.81%  "
$ws.Range("E18").Value = "  +0.79%  "
$ws.Range("E19").Value = "  +0.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "27.240.14"
$ws.Range("E20").Value = "  +0.39%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.093"
$ws.Range("E21").Value = "  +0.65%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.153.99"
$ws.Range("E22").Value = "  +1.68%  "
$ws.Range("E23").Value = "  +0.71%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.440"
$ws.Range("E24").Value = "  +0.89%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.340"
$ws.Range("E25").Value = "  +1.90%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "147.03"
$ws.Range("E26").Value = "  +0.13%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.755"
$ws.Range("E27").Value = "  -0.52%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.25"
$ws.Range("E28").Value = "  +1.02%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "115.20"
$ws.Range("E29").Value = "  +0.69%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.822"
$ws.Range("E30").Value = "  +0.43%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.903"
$ws.Range("E31").Value = "  -0.40%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09287"
$ws.Range("E32").Value = "  +0.91%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05068"
$ws.Range("E33").Value = "  +0.43%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7957"
$ws.Range("E34").Value = "  +3.39%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.238"
$ws.Range("E35").Value = "  +0.34%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.439"
$ws.Range("E36").Value = "  +4.58%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.950"
$ws.Range("E37").Value = "  -1.34%  "
$ws.Range("B38").Value = "TheSandbox"
$ws.Range("C38").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5789"
$ws.Range("E38").Value = "  +3.36%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.583"
$ws.Range("E39").Value = "  -0.62%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02007"
$ws.Range("E40").Value = "  +0.87%  "
$ws.Range("E41").Value = "  +0.16%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.057"
$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.614"
$ws.Range("E43").Value = "  -0.52%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "116.33"
$ws.Range("E44").Value = "  -1.81%  "
$ws.Range("E45").Value = "  +0.96%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4867"
$ws.Range("E46").Value = "  +1.12%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.19"
$ws.Range("E47").Value = "  -0.19%  "
$ws.Range("E48").Value = "  +0.12%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.633"
$ws.Range("E49").Value = "  +1.99%  "
$ws.Range("E50").Value = "  +0.29%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "64.01"
$ws.Range("E51").Value = "  -0.12%  "
